$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "334×2=668" "638×5=3190"
Replace-Text "947×8=7576" "114×6=684"
Replace-Text "289×3=867" "687×9=6183"
Replace-Text "206×5=1030" "534×9=4806"
Replace-Text "491×9=4419" "629×2=1258"
Replace-Text "314×5=1570" "797×3=2391"
Replace-Text "847×3=2541" "119×3=357"
Replace-Text "128×5=640" "146×7=1022"
Replace-Text "914×5=4570" "493×2=986"
Replace-Text "125×3=375" "301×2=602"
Replace-Text "831×6=4986" "989×3=2967"
Replace-Text "650×8=5200" "743×4=2972"
Replace-Text "476×4=1904" "397×2=794"
Replace-Text "847×2=1694" "926×4=3704"
Replace-Text "826×4=3304" "891×6=5346"
Replace-Text "658×9=5922" "919×5=4595"
Replace-Text "479×4=1916" "515×6=3090"
Replace-Text "899×4=3596" "932×6=5592"
Replace-Text "646×9=5814" "180×3=540"
Replace-Text "885×6=5310" "612×2=1224"
Replace-Text "877×5=4385" "996×4=3984"
Replace-Text "985×4=3940" "806×2=1612"
Replace-Text "616×3=1848" "560×9=5040"
Replace-Text "627×9=5643" "102×2=204"
Replace-Text "160×3=480" "518×5=2590"
